# Update the scraped Airbnb experience URLs in column A with new
# federatedSearchId / searchId / sectionId query-string values.
# (Row 1 is the "Link" header; rows 2-21 hold the URLs.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldFederatedSearchId = "526f5db2-2cae-4623-ac66-a5a254ae550c"
$oldSearchId           = "10e5ef69-49a1-4355-8dd6-d4fbd63addef"
$oldSectionId           = "8d95b00f-bf58-4e5e-bb2a-cd27e7f450cc"

$newFederatedSearchId = "e74343bc-d4f2-44b5-bcbb-261b1bd0e7d3"
$newSearchId           = "b06b2908-16d1-40ae-a1bb-0864f164d27b"
$newSectionId           = "7f116f53-1a8d-42d9-a67a-0055d302498a"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $value = $cell.Value()
    if ($value -ne $null) {
        $newValue = $value.Replace($oldFederatedSearchId, $newFederatedSearchId).Replace($oldSearchId, $newSearchId).Replace($oldSectionId, $newSectionId)
        if ($newValue -ne $value) {
            $cell.Value = $newValue
        }
    }
}
